# Auto-generated Excel COM-interop edit script
# Applies the weekly CompStat 83rd Precinct data refresh described in the diff:
#  - volume/date header text bumped to the next week
#  - crime-statistics table (rows 14-29) updated with new counts/percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: "Volume 30   Number  48" -> "...49" ---
$ws.Range("A8").Value = "Volume 30   Number  49"

# --- Report week text: dates roll forward one week ---
$ws.Range("C9").Value = "Report Covering the Week  12/4/2023  Through  12/10/2023"

# --- Cells that change shape (number <-> text placeholder) ---
# Text placeholders reuse the workbook's existing "0" / "***.*" shared strings
# by copying an already-correctly-styled template cell, exactly as Excel's own
# fill/copy would, so style indices match cells like D14/E14 elsewhere in the sheet.
$ws.Range("D14").Copy($ws.Range("C14"))
$ws.Range("D14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("D14").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))
$ws.Range("D14").Copy($ws.Range("F23"))
$ws.Range("D14").Copy($ws.Range("C27"))
$ws.Range("F14").Copy($ws.Range("D27"))
$ws.Range("D27").Value2 = 1
$ws.Range("K14").Copy($ws.Range("E27"))
$ws.Range("E27").Value2 = -100

# --- Plain numeric refreshes (counts + computed %, same cell styling) ---
# Row 15
$ws.Range("I15").Value2 = 22
$ws.Range("K15").Value2 = -42.105263157894
$ws.Range("L15").Value2 = -21.428571428571
$ws.Range("M15").Value2 = -33.333333333333
$ws.Range("N15").Value2 = -63.934426229508
# Row 16
$ws.Range("C16").Value2 = 5
$ws.Range("D16").Value2 = 6
$ws.Range("E16").Value2 = -16.666666666666
$ws.Range("F16").Value2 = 29
$ws.Range("G16").Value2 = 25
$ws.Range("H16").Value2 = 16
$ws.Range("I16").Value2 = 366
$ws.Range("J16").Value2 = 308
$ws.Range("K16").Value2 = 18.831168831168
$ws.Range("L16").Value2 = 56.410256410256
$ws.Range("M16").Value2 = -13.270142180094
$ws.Range("N16").Value2 = -74.776016540317
# Row 17
$ws.Range("C17").Value2 = 10
$ws.Range("D17").Value2 = 8
$ws.Range("E17").Value2 = 25
$ws.Range("F17").Value2 = 34
$ws.Range("G17").Value2 = 31
$ws.Range("H17").Value2 = 9.677419354838
$ws.Range("I17").Value2 = 449
$ws.Range("J17").Value2 = 360
$ws.Range("K17").Value2 = 24.722222222222
$ws.Range("L17").Value2 = 47.213114754098
$ws.Range("M17").Value2 = 25.069637883008
$ws.Range("N17").Value2 = -44.014962593516
# Row 18
$ws.Range("C18").Value2 = 8
$ws.Range("D18").Value2 = 8
$ws.Range("E18").Value2 = 0
$ws.Range("F18").Value2 = 26
$ws.Range("G18").Value2 = 23
$ws.Range("H18").Value2 = 13.043478260869
$ws.Range("I18").Value2 = 237
$ws.Range("J18").Value2 = 288
$ws.Range("K18").Value2 = -17.708333333333
$ws.Range("L18").Value2 = -6.692913385826
$ws.Range("M18").Value2 = -46.380090497737
$ws.Range("N18").Value2 = -81.527669524551
# Row 19
$ws.Range("C19").Value2 = 14
$ws.Range("D19").Value2 = 9
$ws.Range("E19").Value2 = 55.555555555555
$ws.Range("F19").Value2 = 44
$ws.Range("G19").Value2 = 46
$ws.Range("H19").Value2 = -4.347826086956
$ws.Range("I19").Value2 = 676
$ws.Range("J19").Value2 = 666
$ws.Range("K19").Value2 = 1.501501501501
$ws.Range("L19").Value2 = 26.355140186915
$ws.Range("M19").Value2 = 133.103448275862
$ws.Range("N19").Value2 = 26.119402985074
# Row 20
$ws.Range("C20").Value2 = 5
$ws.Range("D20").Value2 = 4
$ws.Range("E20").Value2 = 25
$ws.Range("F20").Value2 = 23
$ws.Range("G20").Value2 = 13
$ws.Range("H20").Value2 = 76.923076923076
$ws.Range("I20").Value2 = 212
$ws.Range("J20").Value2 = 231
$ws.Range("K20").Value2 = -8.225108225108
$ws.Range("L20").Value2 = 27.710843373494
$ws.Range("M20").Value2 = 30.864197530864
$ws.Range("N20").Value2 = -76.548672566371
# Row 21
$ws.Range("C21").Value2 = 42
$ws.Range("D21").Value2 = 35
$ws.Range("E21").Value2 = 20
$ws.Range("F21").Value2 = 158
$ws.Range("G21").Value2 = 139
$ws.Range("H21").Value2 = 13.669064748201
$ws.Range("I21").Value2 = 1966
$ws.Range("J21").Value2 = 1894
$ws.Range("K21").Value2 = 3.801478352692
$ws.Range("L21").Value2 = 28.664921465968
$ws.Range("M21").Value2 = 14.235909355026
$ws.Range("N21").Value2 = -61.276344297813
# Row 22
$ws.Range("G22").Value2 = 3
$ws.Range("H22").Value2 = -33.333333333333
$ws.Range("M22").Value2 = -7.692307692307
# Row 23
$ws.Range("H23").Value2 = -100
$ws.Range("L23").Value2 = 12.121212121212
$ws.Range("M23").Value2 = 60.869565217391
# Row 24
$ws.Range("C24").Value2 = 16
$ws.Range("D24").Value2 = 16
$ws.Range("E24").Value2 = 0
$ws.Range("F24").Value2 = 81
$ws.Range("G24").Value2 = 72
$ws.Range("H24").Value2 = 12.5
$ws.Range("I24").Value2 = 933
$ws.Range("J24").Value2 = 997
$ws.Range("K24").Value2 = -6.419257773319
$ws.Range("L24").Value2 = 1.413043478260
$ws.Range("M24").Value2 = 28.867403314917
# Row 25
$ws.Range("C25").Value2 = 7
$ws.Range("D25").Value2 = 15
$ws.Range("E25").Value2 = -53.333333333333
$ws.Range("F25").Value2 = 47
$ws.Range("G25").Value2 = 37
$ws.Range("H25").Value2 = 27.027027027027
$ws.Range("I25").Value2 = 646
$ws.Range("J25").Value2 = 555
$ws.Range("K25").Value2 = 16.396396396396
$ws.Range("L25").Value2 = 29.718875502008
$ws.Range("M25").Value2 = -13.751668891855
# Row 26
$ws.Range("I26").Value2 = 39
$ws.Range("K26").Value2 = -25
$ws.Range("L26").Value2 = -9.302325581395
# Row 27
$ws.Range("F27").Value2 = 2
$ws.Range("H27").Value2 = -33.333333333333
$ws.Range("I27").Value2 = 59
$ws.Range("J27").Value2 = 62
$ws.Range("K27").Value2 = -4.838709677419
$ws.Range("L27").Value2 = 20.408163265306
# Row 28
$ws.Range("D28").Value2 = 2
$ws.Range("G28").Value2 = 3
$ws.Range("H28").Value2 = -33.333333333333
$ws.Range("J28").Value2 = 17
$ws.Range("K28").Value2 = -11.764705882352
$ws.Range("L28").Value2 = -60.526315789473
$ws.Range("N28").Value2 = -93.181818181818
# Row 29
$ws.Range("G29").Value2 = 2
$ws.Range("H29").Value2 = 0
$ws.Range("J29").Value2 = 16
$ws.Range("K29").Value2 = -12.5
$ws.Range("L29").Value2 = -50
$ws.Range("N29").Value2 = -93.103448275862
